$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Row 8 previously held a single (empty, thick-bottom-bordered) cell left
# over from the table's closing border. Delete it outright so the new test
# case row starts from a clean slate instead of inheriting that row-level
# "thick bottom border" flag.
$ws.Rows.Item(8).Delete()

# --- New test case row (row 8). Shared strings must be created in this
# --- exact order so they land at the same indices as the target workbook.
$ws.Range("C8").Value = "Verifies that there are only two type of account allowed"
$ws.Range("F8").Value = "1. Call the createAcc account parameters for accountType, owner respectively"
$ws.Range("H8").Value = "1. Throws an exception saying that invalid type of account provided "
$ws.Range("I8").Value = "1. Throws an exception which states, invald type of account provided"
$ws.Range("K8").Value = "CreateAccount class method"
$ws.Range("B8").Value = "TC007"
$ws.Range("G8").Value = "1. Accounttype: student                  2. Accountname: John"

$ws.Range("A8").Value = "TS01"
$ws.Range("D8").Value = "N/A"
$ws.Range("E8").Value = "N/A"
$ws.Range("J8").Value = "Pass"
$ws.Range("L8").Value = "Sanjay Sohal"
$ws.Range("N8").Value = "Sanjay Sohal"
$ws.Range("P8").Value = "OS: Windows 8.1                   IDE: Eclipse"

# Copy formatting from row 7 (wrap-text body style / date style) onto row 8
# so the new row reuses the existing style indices instead of minting new ones.
$ws.Range("A7:L7").Copy()
$ws.Range("A8:L8").PasteSpecial(-4122)
$ws.Range("N7").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("P7").Copy()
$ws.Range("P8").PasteSpecial(-4122)

$ws.Range("M7:O7").Copy()
$ws.Range("M8:O8").PasteSpecial(-4122)
$ws.Range("M8").Value = 42107
$ws.Range("O8").Value = 42107

$excel.CutCopyMode = $false

$ws.Rows.Item(8).RowHeight = 84.75

$ws.Activate()
$ws.Range("M8").Select()

Write-Host "done"
